# IronCalc COM-interop edit script
# Reflects a refresh of the COVID-19 "paises" dataset on sheet "Pais":
#   - the "last updated" timestamp in A1 moves from 22:26 to 23:43
#   - several countries receive new totals, which shifts their rank in the
#     (descending, by "Casos totales") list, so some rows below a re-ranked
#     country keep their old numbers but show the country that is now in that
#     rank position (e.g. Costa Rica jumps above Armenia/Kirguistan; Suazilandia
#     above Nicaragua; Tunez above Ruanda/Cuba/Mozambique; Islas Malvinas above Montserrat).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Refresh the "datos actualizados" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 3 de Septiembre de 2020 a las 23:43"

# 2. Apply the updated figures / re-ranked country labels, row by row
# Row 4
$ws.Range("B4").Value = 6328844
$ws.Range("C4").Value = 38107
$ws.Range("D4").Value = 3566046
$ws.Range("E4").Value = 2571888
$ws.Range("G4").Value = 946
$ws.Range("H4").Value = 190910

# Row 5
$ws.Range("B5").Value = 4040163
$ws.Range("C5").Value = 38741
$ws.Range("E5").Value = 705107
$ws.Range("G5").Value = 752
$ws.Range("H5").Value = 124651

# Row 23
$ws.Range("B23").Value = 248814
$ws.Range("C23").Value = 1423
$ws.Range("E23").Value = 16315
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = 9399

# Row 29
$ws.Range("B29").Value = 124455
$ws.Range("C29").Value = 2991
$ws.Range("D29").Value = 98645
$ws.Range("E29").Value = 24825

# Row 59 -> Costa Rica
$ws.Range("A59").Value = "Costa Rica"
$ws.Range("B59").Value = 44458
$ws.Range("C59").Value = 1153
$ws.Range("D59").Value = 17855
$ws.Range("E59").Value = 26143
$ws.Range("G59").Value = 7
$ws.Range("H59").Value = 460

# Row 60 -> Armenia
$ws.Range("A60").Value = "Armenia"
$ws.Range("B60").Value = 44271
$ws.Range("C60").Value = 196
$ws.Range("D60").Value = 38855
$ws.Range("E60").Value = 4529
$ws.Range("G60").Value = 3
$ws.Range("H60").Value = 887

# Row 61 -> Kirguistan
$ws.Range("A61").Value = "Kirguistan"
$ws.Range("B61").Value = 44135
$ws.Range("C61").Value = 99
$ws.Range("D61").Value = 39174
$ws.Range("E61").Value = 3902
$ws.Range("H61").Value = 1059

# Row 81
$ws.Range("B81").Value = 18208
$ws.Range("C81").Value = 47
$ws.Range("D81").Value = 17045
$ws.Range("E81").Value = 1044
$ws.Range("G81").Value = 2
$ws.Range("H81").Value = 119

# Row 83
$ws.Range("B83").Value = 16775
$ws.Range("C83").Value = 158
$ws.Range("D83").Value = 11935
$ws.Range("E83").Value = 4182
$ws.Range("G83").Value = 10
$ws.Range("H83").Value = 658

# Row 104
$ws.Range("B104").Value = 7106
$ws.Range("C104").Value = 17
$ws.Range("D104").Value = 6588
$ws.Range("E104").Value = 358

# Row 107
$ws.Range("B107").Value = 6678
$ws.Range("C107").Value = 40
$ws.Range("D107").Value = 5263
$ws.Range("E107").Value = 1209

# Row 108
$ws.Range("B108").Value = 5593
$ws.Range("C108").Value = 14
$ws.Range("D108").Value = 3516
$ws.Range("E108").Value = 1902

# Row 113
$ws.Range("B113").Value = 4729
$ws.Range("C113").Value = 17
$ws.Range("D113").Value = 1807
$ws.Range("E113").Value = 2860

# Row 114 -> Suazilandia
$ws.Range("A114").Value = "Suazilandia"
$ws.Range("B114").Value = 4720
$ws.Range("C114").Value = 52
$ws.Range("D114").Value = 3789
$ws.Range("E114").Value = 837
$ws.Range("H114").Value = 94

# Row 115 -> Nicaragua
$ws.Range("A115").Value = "Nicaragua"
$ws.Range("D115").Value = 2913
$ws.Range("E115").Value = 1614
$ws.Range("H115").Value = 141

# Row 117 -> Tunez
$ws.Range("A117").Value = "Tunez"
$ws.Range("B117").Value = 4394
$ws.Range("C117").Value = 198
$ws.Range("D117").Value = 1681
$ws.Range("E117").Value = 2629
$ws.Range("G117").Value = 3
$ws.Range("H117").Value = 84

# Row 118 -> Ruanda
$ws.Range("A118").Value = "Ruanda"
$ws.Range("B118").Value = 4255
$ws.Range("C118").Value = 37
$ws.Range("D118").Value = 2163
$ws.Range("E118").Value = 2074
$ws.Range("G118").Value = 1
$ws.Range("H118").Value = 18

# Row 119 -> Cuba
$ws.Range("A119").Value = "Cuba"
$ws.Range("B119").Value = 4214
$ws.Range("C119").Value = 88
$ws.Range("D119").Value = 3474
$ws.Range("E119").Value = 640
$ws.Range("G119").Value = 2
$ws.Range("H119").Value = 100

# Row 120 -> Mozambique
$ws.Range("A120").Value = "Mozambique"
$ws.Range("B120").Value = 4207
$ws.Range("C120").Value = 90
$ws.Range("D120").Value = 2370
$ws.Range("E120").Value = 1811
$ws.Range("G120").Value = 1
$ws.Range("H120").Value = 26

# Row 154
$ws.Range("B154").Value = 1443
$ws.Range("C154").Value = 9
$ws.Range("D154").Value = 1071
$ws.Range("E154").Value = 341
$ws.Range("G154").Value = 1
$ws.Range("H154").Value = 31

# Row 157
$ws.Range("B157").Value = 1375
$ws.Range("C157").Value = 5
$ws.Range("E157").Value = 245

# Row 189
$ws.Range("B189").Value = 177
$ws.Range("C189").Value = 1
$ws.Range("D189").Value = 150

# Row 214 -> Islas Malvinas
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

# Row 215 -> Montserrat
$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1

Write-Host "Updated paises sheet: timestamp + 24 data rows"
